$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2022-01-18T12:41:26-05:00"

# Update the "epi-questions" concept row on the Concepts sheet
$conceptsWs = $wb.Worksheets.Item("Concepts")
$conceptsWs.Range("B28").Value = "epi-observations"
$conceptsWs.Range("C28").Value = "Epi Observations"
